# Customer login workflow implemented and tested
# End to end workflow (Bank manager & Customer) implemented and tested
#
# Adds a Deposit / Withdrawl section (columns J/K/L) to the
# "PostitiveEntries" sheet, mirroring the existing ${FNAME}/${LNAME}/...
# header style, and switches the active sheet/selection to reflect the
# author's final view state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # PostitiveEntries
$ws2 = $wb.Worksheets.Item(2)   # NegativeEntries

# --- New header cells on PostitiveEntries row 1 (J1:L1), matching the
#     existing B1/D1/F1/H1 yellow-fill header style ---
$ws1.Range("J1").Value = '${DEPOSIT}'
$ws1.Range("J1").Interior.Color = 65535

$ws1.Range("K1").Interior.Color = 65535

$ws1.Range("L1").Value = '${WITHDRAWL}'
$ws1.Range("L1").Interior.Color = 65535

# --- New data: Deposit (J) / Withdrawl (L) amounts for rows 2-5 ---
$ws1.Range("J2").Value = 1000
$ws1.Range("L2").Value = 500

$ws1.Range("J3").Value = 1234
$ws1.Range("L3").Value = 123

$ws1.Range("J4").Value = 89500
$ws1.Range("L4").Value = 84000

$ws1.Range("J5").Value = 1100
$ws1.Range("L5").Value = 5

# --- Column widths for the new columns ---
$ws1.Columns.Item(10).ColumnWidth = 9.5               # J
$ws1.Columns.Item(12).ColumnWidth = 13.333333333333334 # L
$ws2.Columns.Item(10).ColumnWidth = 9.5               # J (width-only on NegativeEntries)

# --- Final view state: NegativeEntries selection moves to J1, and
#     PostitiveEntries becomes the active sheet with L10 selected ---
[void]$ws2.Range("J1").Select()

[void]$ws1.Activate()
[void]$ws1.Range("L10").Select()
